$d = $word.ActiveDocument

# 1. Merge the split "I" + "nXtremis..." runs at the start of the document
#    into a single contiguous run by replacing "InXtremis" with itself via
#    paragraph-scoped Find (forces run consolidation across the split).
$p1 = $d.Paragraphs(1).Range
$p1.Find.Execute("InXtremis", $false, $false, $false, $false, $false, `
                  $true, 1, $false, "InXtremis", 2) | Out-Null

# 2. Merge the "Interface, collecte de données, modularité, libre d'accès"
#    paragraph with the following paragraph's text, splitting the combined
#    text after "tout ça dans un ".
$apostrophe = [char]0x2019

$oldText = "Interface, collecte de données, modularité, libre d" + $apostrophe + "accès"
$newText = "Il y a beaucoup de « mallettes d" + $apostrophe + "évasion » sur le marché, mais 4 éléments nous distinguent des compétiteurs. L" + $apostrophe + "interface graphique que procure l" + $apostrophe + "écran, la collecte de données pour une meilleure gestion, la modularité des énigmes et le libre d" + $apostrophe + "accès de la mallette, tout ça dans un "

$secondOld = "Il y a beaucoup de « mallettes d" + $apostrophe + "évasion » sur le marché, mais 4 éléments nous distinguent des compétiteurs. L" + $apostrophe + "interface graphique que procure l" + $apostrophe + "écran, la collecte de données pour une meilleure gestion, la modularité des énigmes et le libre d" + $apostrophe + "accès de la mallette, tout ça dans un emballage"
$secondNew = "emballage"

# Find the paragraph that currently reads "Interface, collecte de données, ..."
# and replace it with the first half of the following paragraph's text.
$found = $false
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs($i)
    if ($para.Range.Text.TrimEnd([char]0x0D, [char]0x07) -eq $oldText) {
        $nextPara = $d.Paragraphs($i + 1)
        $nextRange = $nextPara.Range

        # Replace the "Interface, ..." paragraph text with the new first part.
        $paraRange = $para.Range
        $paraRange.MoveEnd(1, -1) | Out-Null  # exclude paragraph mark
        $paraRange.Text = $newText

        $found = $true
        break
    }
}

# Trim the following paragraph so it starts at "emballage compact et durable..."
$d.Content.Find.Execute($secondOld, $false, $false, $false, $false, $false, `
                         $true, 1, $false, $secondNew, 2) | Out-Null

# Now remove the paragraph break between the merged "Interface" paragraph
# (which now starts with "Il y a beaucoup...") and the "emballage..." paragraph,
# so they become a single paragraph with two runs.
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs($i)
    if ($para.Range.Text.StartsWith("Il y a beaucoup de")) {
        $endRange = $d.Range($para.Range.End - 1, $para.Range.End)
        $endRange.Text = ""
        break
    }
}
